$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks (G2:G4) first so we can cleanly re-add all
# five (G2:G6) in order further down -- re-adding after deleting renumbers
# the hyperlink relationship ids starting at rId1 again.
$ws.Range("G2").Hyperlinks.Delete()

# --- Column A ---
$ws.Range("A2").Value = "IBBI"
$ws.Range("A3").Value = "IBBI"
$ws.Range("A4").Value = "IBBI"
$ws.Range("A5").Value = "IBBI"
$ws.Range("A6").Value = "IBBI"

# --- Column B ---
$ws.Range("B2").Value = "Acts"
$ws.Range("B3").Value = "Discussion Paper"
$ws.Range("B4").Value = "Discussion Paper"
$ws.Range("B5").Value = "Guidelines"
$ws.Range("B6").Value = "Guidelines"

# --- Column C ---
$ws.Range("C2").Value = "2025"
$ws.Range("C3").Value = "2025"
$ws.Range("C4").Value = "2025"
$ws.Range("C5").Value = "2025"
$ws.Range("C6").Value = "2025"

# --- Column D ---
$ws.Range("D2").Value = "August"
$ws.Range("D3").Value = "August"
$ws.Range("D4").Value = "August"
$ws.Range("D5").Value = "August"
$ws.Range("D6").Value = "August"

# --- Column E ---
$ws.Range("E2").Value = "2025-08-12"
$ws.Range("E3").Value = "2025-08-12"
$ws.Range("E4").Value = "2025-08-12"
$ws.Range("E5").Value = "2025-08-13"
$ws.Range("E6").Value = "2025-08-13"

# --- Column F ---
$ws.Range("F2").Value = "Bill - The Insolvency and Bankruptcy Code (Amendment) Bill, 2025 (as introduced in Lok Sabha)"
$ws.Range("F3").Value = "Discussion paper - Review of Limit on Number of Assignments by IPs"
$ws.Range("F4").Value = "Discussion paper on deletion of Clause 6 from the Code of Conduct for Insolvency Professionals"
$ws.Range("F5").Value = "Amendments to the Insolvency and Bankruptcy Board of India (Continuing Professional Education for Insolvency Professionals) Guidelines, 2019 (220.55 KB)"
$ws.Range("F6").Value = "Insolvency and Bankruptcy Board of India (Continuing Professional Education for Insolvency Professionals) Guidelines, 2019 (Updated as on 13th August, 2025) (302.56 KB)"

# --- Column G (PDF URLs, gets the Hyperlink style + real hyperlinks below) ---
$ws.Range("G2").Value = "https://ibbi.gov.in/uploads/legalframwork/da78600a457741799bb2e7c8da25f946.pdf"
$ws.Range("G3").Value = "https://ibbi.gov.in/uploads/public_comments/Discussion paper - Review of Limit on Number of Assignments by IPs - final.pdf"
$ws.Range("G4").Value = "https://ibbi.gov.in/uploads/public_comments/Discussion paper on deletion of Clause 6 from the Code of Conduct for Insolvency Professionals - final.pdf"
$ws.Range("G5").Value = "https://ibbi.gov.in/uploads/legalframwork/2025-08-13-171527-5bcdm-290d5d85373e4dec2c48e209925a1bbb.pdf"
$ws.Range("G6").Value = "https://ibbi.gov.in/uploads/legalframwork/69e8b480fe681423a04b7a93b34ecd31.pdf"

# --- Column H ---
$ws.Range("H2").Value = "Bill_The_Insolvency_and_Bankruptcy_Code_Amendment_Bill_2025_as_introduced_in_Lok_ad048024.pdf"
$ws.Range("H3").Value = "Discussion_paper_Review_of_Limit_on_Number_of_Assignments_by_IPs_f73af763.pdf"
$ws.Range("H4").Value = "Discussion_paper_on_deletion_of_Clause_6_from_the_Code_of_Conduct_for_Insolvency_68137bd6.pdf"
$ws.Range("H5").Value = "Amendments_to_the_Insolvency_and_Bankruptcy_Board_of_India_Continuing_Profession_c9bd3f17.pdf"
$ws.Range("H6").Value = "Insolvency_and_Bankruptcy_Board_of_India_Continuing_Professional_Education_for_I_7c21a393.pdf"

# --- Column I ---
$ws.Range("I2").Value = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/IBBI/Acts/2025/August/Bill_The_Insolvency_and_Bankruptcy_Code_Amendment_Bill_2025_as_introduced_in_Lok_ad048024.pdf"
$ws.Range("I3").Value = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/IBBI/Discussion Paper/2025/August/Discussion_paper_Review_of_Limit_on_Number_of_Assignments_by_IPs_f73af763.pdf"
$ws.Range("I4").Value = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/IBBI/Discussion Paper/2025/August/Discussion_paper_on_deletion_of_Clause_6_from_the_Code_of_Conduct_for_Insolvency_68137bd6.pdf"
$ws.Range("I5").Value = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/IBBI/Guidelines/2025/August/Amendments_to_the_Insolvency_and_Bankruptcy_Board_of_India_Continuing_Profession_c9bd3f17.pdf"
$ws.Range("I6").Value = "/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/IBBI/Guidelines/2025/August/Insolvency_and_Bankruptcy_Board_of_India_Continuing_Professional_Education_for_I_7c21a393.pdf"

# Re-create the hyperlinks on G2:G6 pointing at the new PDF URLs.
$ws.Hyperlinks.Add($ws.Range("G2"), "https://ibbi.gov.in/uploads/legalframwork/da78600a457741799bb2e7c8da25f946.pdf")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://ibbi.gov.in/uploads/public_comments/Discussion paper - Review of Limit on Number of Assignments by IPs - final.pdf")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://ibbi.gov.in/uploads/public_comments/Discussion paper on deletion of Clause 6 from the Code of Conduct for Insolvency Professionals - final.pdf")
$ws.Hyperlinks.Add($ws.Range("G5"), "https://ibbi.gov.in/uploads/legalframwork/2025-08-13-171527-5bcdm-290d5d85373e4dec2c48e209925a1bbb.pdf")
$ws.Hyperlinks.Add($ws.Range("G6"), "https://ibbi.gov.in/uploads/legalframwork/69e8b480fe681423a04b7a93b34ecd31.pdf")

# Make sure G2:G6 keep the workbook's existing "Hyperlink" cell style (rather
# than whatever extra style index Hyperlinks.Add allocates).
$ws.Range("G2:G6").Style = "Hyperlink"
